$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns N, O, P (with same style as existing header cells)
$ws.Range("N1").Value = "nr of simulations"
$ws.Range("O1").Value = "provided beta"
$ws.Range("P1").Value = "calculated beta"
$ws.Range("M1").Copy()
$ws.Range("N1:P1").PasteSpecial(-4122)  # xlPasteFormats

# Row 2 (customer type A)
$ws.Range("F2").Value = -1.627534966661221
$ws.Range("G2").Value = 811
$ws.Range("H2").Value = 1.188040320916699
$ws.Range("I2").Value = -1.369932432432432
$ws.Range("J2").Value = -287.6858108108108
$ws.Range("K2").Value = -341.7823429988565
$ws.Range("L2").Value = 249.4884673925068
$ws.Range("M2").Value = -0.8707312071637534

# Row 3 (customer type B)
$ws.Range("F3").Value = -1.761145907544019
$ws.Range("G3").Value = 674
$ws.Range("H3").Value = 1.54688186537991
$ws.Range("I3").Value = -1.138513513513514
$ws.Range("J3").Value = -239.0878378378379
$ws.Range("K3").Value = -369.8406405842441
$ws.Range("L3").Value = 324.8451917297811
$ws.Range("M3").Value = -0.9422130605360505

# Row 4 (customer type C)
$ws.Range("F4").Value = -1.068375785164082
$ws.Range("G4").Value = 507
$ws.Range("H4").Value = 1.247492041059441
$ws.Range("I4").Value = -0.856418918918919
$ws.Range("J4").Value = -179.847972972973
$ws.Range("K4").Value = -224.3589148844573
$ws.Range("L4").Value = 261.9733286224827
$ws.Range("M4").Value = -0.5715810450627841

# Row 5 (TOTAL) plus new convergence-criteria columns
$ws.Range("J5").Value = -1.121621621621622
$ws.Range("K5").Value = -1.485685553123108
$ws.Range("L5").Value = 1.324587119651927
$ws.Range("M5").Value = -2.384525312762588
$ws.Range("N5").Value = 8
$ws.Range("O5").Value = 0.02
$ws.Range("P5").Value = 0.001867481408371164
